$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-20 down to 3-21),
# pushing in a new IPO subscription entry for 오상헬스케어 (NH).
$ws.Rows.Item(2).Insert()

# The three date-like text columns need to stay as plain text (matching the
# rest of the sheet), not get auto-converted to Excel date serials.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "2024-03-04"
$ws.Range("B2").Value = "오상헬스케어"
$ws.Range("C2").Value = "NH"
$ws.Range("D2").Value = "2024-03-07"
$ws.Range("E2").Value = "2024-03-13"
$ws.Range("F2").Value = 19800000
$ws.Range("G2").Value = 990000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 13000
$ws.Range("J2").Value = 15000
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "2126.13 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"

# The freshly-inserted row inherits formatting from the row above (header);
# strip it back to the plain (unstyled) look used by the rest of the data
# rows so the new row matches its neighbours.
$ws.Range("A2:T2").ClearFormats()
